$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = -6
    4  = 5
    5  = 3
    6  = 1
    7  = -2
    8  = 3
    9  = -1
    10 = 5
    11 = 8
    12 = 4
    13 = -5
    14 = -1
    16 = 2
    17 = -3
    18 = -3
    19 = -1
    21 = -4
    22 = 1
    23 = -3
    24 = 1
    25 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
